# Generate Report for Handoff
# Updates the localization-status report with the results of the latest
# handoff run: new GUID-named source/target files and refreshed timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "157c079d-e566-4677-b603-711fd8ce4d3d"
$newGuid = "a02157ab-1c54-43ee-9f4e-a990972644b9"
$oldHash = "b36f518f22265336406e5412e346be3e55d1d800"
$newHash = "58c52a7a14f9ae85b7c96ea78a3a668252930f80"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newOverviewDate = "2016-03-24 15:20:54"
$newZhDate        = "2016-03-24 15:20:50"
$newDeDate        = "2016-03-24 15:20:54"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/6a73d23b55a7e87102be477660dbf15a917ce944/e2e/$oldGuid.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36e342265e33aac37c0dcf0d469d962dd032cf58/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/51ad9f1a9fa751edbc8eb06a35f655de16d72a8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, $newMdName)

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName
$wsZh.Range("E2").Value = $newZhDate

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddress, [Type]::Missing, [Type]::Missing, $newZhXlfName)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName
$wsDe.Range("E2").Value = $newDeDate

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddress, [Type]::Missing, [Type]::Missing, $newDeXlfName)
